$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.499.55"
$ws.Range("E2").Value = "  +4.61%  "
$ws.Range("D3").Value = "3.361.95"
$ws.Range("E3").Value = "  +9.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("E7").Value = "  +10.24%  "
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "3.357.23"
$ws.Range("E10").Value = "  +9.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "98.047.19"
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000247"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("E16").Value = "  +9.66%  "
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "3.363.73"
$ws.Range("E18").Value = "  +9.38%  "
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("E23").Value = "  +9.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "3.543.87"
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.257"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.56%  "
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  +3.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "522.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  +4.32%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.450"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("E44").Value = "  +5.69%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("E48").Value = "  +5.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.45%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.25%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.99%  "
